# The worksheet data (name, task, start_date, end_date plus a leading
# id column A) was re-sorted so that it is ordered by the id column (A)
# ascending, instead of the previous ad-hoc ordering by task description.
# Row 1 holds the column headers and must stay in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

$sortRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$key1 = $ws.Cells.Item(1, 1)

$sortRange.Sort($key1, 1, $null, $null, 1, $null, 1, 1)
